$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force D2:E51 to text format so numeric-looking strings are preserved as text
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.386.85'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '1.834.77'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.88%  '
$ws.Range('D5').Value = '314.44'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').Value = '1.010'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').Value = '0.4739'
$ws.Range('E7').Value = '  +1.99%  '
$ws.Range('D8').Value = '0.3693'
$ws.Range('E8').Value = '  +0.99%  '
$ws.Range('D9').Value = '0.07468'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Value = '0.8867'
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('D12').Value = '1.896.77'
$ws.Range('E12').Value = '  +6.83%  '
$ws.Range('D13').Value = '0.07350'
$ws.Range('E13').Value = '  +3.55%  '
$ws.Range('D14').Value = '5.455'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').Value = '93.16'
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('D16').Value = '6.585'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '0.000008824'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = '1.010'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = '27.661.24'
$ws.Range('E20').Value = '  +2.79%  '
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').Value = '2.105.01'
$ws.Range('E24').Value = '  +3.58%  '
$ws.Range('D25').Value = '1.897'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '152.02'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').Value = '18.65'
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('D29').Value = '5.247'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').Value = '117.60'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').Value = '0.7560'
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').Value = '1.176'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').Value = '4.551'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').Value = '2.949'
$ws.Range('E35').Value = '  +1.51%  '
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').Value = '0.05350'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('D39').Value = '0.01952'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = '2.978'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').Value = '7.329'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('D42').Value = '2.409'
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('D43').Value = '0.5335'
$ws.Range('D44').Value = '0.1659'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '8.509'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('D47').Value = '10.54'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').Value = '105.03'
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').Value = '1.676'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('D51').Value = '0.06293'
$ws.Range('E51').Value = '  -0.07%  '

# Restore default (unstyled) formatting so only cell content differs from before
$dataRange.ClearFormats()
